# Organization-Firewall-Manager-Architecture.pptx edit script
# Issue #22 fixes: refresh footer date fields, de-highlight two background
# rectangles, nudge several diagram shapes, recolor the "VPC" label, and
# add two more AWS-logo graphics near the CloudFormation icons.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---------------------------------------------------------------
# 1) Remove the translucent fill from the two big background
#    rectangles (Rectangle 128 / Rectangle 126) -> <a:noFill/>
# ---------------------------------------------------------------
$s.Shapes.Item(1).Fill.Visible = $false   # Rectangle 128
$s.Shapes.Item(4).Fill.Visible = $false   # Rectangle 126

# ---------------------------------------------------------------
# 2) Reposition / resize a batch of diagram shapes.
# ---------------------------------------------------------------
# 3: pic id126 Graphic 2
$sh = $s.Shapes.Item(3)
$sh.Left = 160.83118510236218
$sh.Top = 19.183780527559055

# 5: TextBox 130
$sh = $s.Shapes.Item(5)
$sh.Left = 598.2204904409448
$sh.Top = 130.0444114488189

# 6: TextBox 2 (id140)
$sh = $s.Shapes.Item(6)
$sh.Left = 418.58213898425197
$sh.Top = 137.34110336220473

# 7: Graphic 1
$sh = $s.Shapes.Item(7)
$sh.Left = 454.7314920629921
$sh.Top = 190.5307854015748

# 8: TextBox 2 (id142)
$sh = $s.Shapes.Item(8)
$sh.Left = 415.01313860629926
$sh.Top = 234.24385926771654

# 9: Graphic 142
$sh = $s.Shapes.Item(9)
$sh.Left = 452.93550181102364
$sh.Top = 90.55614173228346

# 10: TextBox 143
$sh = $s.Shapes.Item(10)
$sh.Left = 512.1583564566929
$sh.Top = 174.2552725905512

# 11: Graphic 144
$sh = $s.Shapes.Item(11)
$sh.Left = 533.7808541417323
$sh.Top = 142.6948778897638

# 12: Straight Arrow Connector 147
$sh = $s.Shapes.Item(12)
$sh.Left = 521.9282537165354
$sh.Top = 112.6207854015748

# 13: Straight Arrow Connector 148
$sh = $s.Shapes.Item(13)
$sh.Left = 474.343765527559
$sh.Top = 156.43653143307088

# 14: TextBox 150
$sh = $s.Shapes.Item(14)
$sh.Left = 262.87062092125984
$sh.Top = 177.34810723622047
$sh.Width = 95.25243794488188
$sh.Height = 36.351576803149605

# 15: Graphic 151
$sh = $s.Shapes.Item(15)
$sh.Left = 288.2499242598425
$sh.Top = 133.24109736220473

# 16: Rectangle 152
$sh = $s.Shapes.Item(16)
$sh.Left = 385.3826761653543
$sh.Top = 83.63787501574802
$sh.Width = 319.4263769527559
$sh.Height = 186.78409648818896

# 17: Straight Arrow Connector 153
$sh = $s.Shapes.Item(17)
$sh.Left = 336.3624429448819
$sh.Top = 150.53299812598428

# 24: Oval 225
$sh = $s.Shapes.Item(24)
$sh.Left = 290.15095588188973
$sh.Top = 49.86015948031496

# 25: Oval 226
$sh = $s.Shapes.Item(25)
$sh.Left = 289.23661817322835
$sh.Top = 133.76905111811024

# 26: Oval 227
$sh = $s.Shapes.Item(26)
$sh.Left = 452.93550181102364
$sh.Top = 90.56189376377954

# 27: Graphic 77
$sh = $s.Shapes.Item(27)
$sh.Left = 608.5525216850393
$sh.Top = 88.08023522047245

# 28: Oval 228
$sh = $s.Shapes.Item(28)
$sh.Left = 543.652313464567
$sh.Top = 130.86377752755905

# 29: Oval 229
$sh = $s.Shapes.Item(29)
$sh.Left = 455.93678365354333
$sh.Top = 191.0012598425197

# 30: Oval 230
$sh = $s.Shapes.Item(30)
$sh.Left = 609.5605781811024
$sh.Top = 88.60567129133858

# 31: Oval 231
$sh = $s.Shapes.Item(31)
$sh.Left = 288.1870886141732
$sh.Top = 290.471817023622

# ---------------------------------------------------------------
# 3) Recolor the "VPC" text run from accent5 to a fixed navy srgb.
# ---------------------------------------------------------------
$vpc = $s.Shapes.Item(36)
$vpc.TextFrame.TextRange.Font.Color.RGB = 6299648   # RGB(0x00,0x20,0x60) => 002060

# ---------------------------------------------------------------
# 4) Duplicate the AWS-logo graphic (shape 3, "Graphic 2") twice,
#    and move the copies next to the CloudFormation icons.
# ---------------------------------------------------------------
$newPic1 = $s.Shapes.Item(3).Duplicate()
$newPic1.Left = 258.71653843307087
$newPic1.Top = 45.003072866141736

$newPic2 = $s.Shapes.Item(3).Duplicate()
$newPic2.Left = 257.95054718110237
$newPic2.Top = 285.7607884015748

# ---------------------------------------------------------------
# 5) Refresh the cached "today" footer date on the slide master and
#    every slide layout (the field keeps auto-updating, but the
#    cached <a:t> text is refreshed the way PowerPoint does on save).
# ---------------------------------------------------------------
$newDate = "4/28/21"

$master = $p.SlideMaster
foreach ($shMaster in $master.Shapes) {
    if ($shMaster.Name -like "Date Placeholder*") {
        $shMaster.TextFrame.TextRange.Text = $newDate
    }
}

for ($i = 1; $i -le $master.CustomLayouts.Count; $i++) {
    $layout = $master.CustomLayouts.Item($i)
    foreach ($shLayout in $layout.Shapes) {
        if ($shLayout.Name -like "Date Placeholder*") {
            $shLayout.TextFrame.TextRange.Text = $newDate
        }
    }
}
